$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Changes of 31st March 2022": rows 2-22 get new PackageTrackNum (col C) /
# ShipmentTrackNum (col D) values. Column D mirrors column C only for the
# rows that already carried a value there (5,6,7,13,14,15,16,17).
$rowValues = @(
    @{ Row = 2;  Value = "320018191948"; HasD = $false }
    @{ Row = 3;  Value = "320018191959"; HasD = $false }
    @{ Row = 4;  Value = "320018191981"; HasD = $false }
    @{ Row = 5;  Value = "320018192006"; HasD = $true }
    @{ Row = 6;  Value = "320018192040"; HasD = $true }
    @{ Row = 7;  Value = "320018192061"; HasD = $true }
    @{ Row = 8;  Value = "320018192094"; HasD = $false }
    @{ Row = 9;  Value = "320018192131"; HasD = $false }
    @{ Row = 10; Value = "320018192164"; HasD = $false }
    @{ Row = 11; Value = "320018192186"; HasD = $false }
    @{ Row = 12; Value = "320018192223"; HasD = $false }
    @{ Row = 13; Value = "320018192245"; HasD = $true }
    @{ Row = 14; Value = "320018192278"; HasD = $true }
    @{ Row = 15; Value = "320018192290"; HasD = $true }
    @{ Row = 16; Value = "320018192326"; HasD = $true }
    @{ Row = 17; Value = "320018192348"; HasD = $true }
    @{ Row = 18; Value = "320018192381"; HasD = $false }
    @{ Row = 19; Value = "320018192407"; HasD = $false }
    @{ Row = 20; Value = "320018192430"; HasD = $false }
    @{ Row = 21; Value = "320018192451"; HasD = $false }
    @{ Row = 22; Value = "320018192484"; HasD = $false }
)

foreach ($entry in $rowValues) {
    $row = $entry.Row
    $val = $entry.Value

    # A leading apostrophe forces the numeric-looking string to be stored
    # as text (t="s") instead of being coerced to a number; resetting the
    # style back to "Normal" afterwards drops the transient quote-prefix
    # flag that the apostrophe entry adds, so the cell keeps its original
    # (default) style - matching these cells' pre-edit formatting.
    $cRange = $ws.Range("C$row")
    $cRange.Value = "'" + $val
    $cRange.Style = "Normal"

    if ($entry.HasD) {
        $dRange = $ws.Range("D$row")
        $dRange.Value = "'" + $val
        $dRange.Style = "Normal"
    }
}
